$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $val) {
    # Plain string value; safe for values Excel would not mistake for a number
    $ws.Range($addr).Value = $val
}

function Set-ForcedText($addr, $val) {
    # Force text storage so Excel doesn't auto-convert numeric-looking strings
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - Bitcoin
Set-Text "D2" "70.879.29"
Set-Text "E2" "  -0.22%  "

# Row 3 - Ethereum
Set-Text "D3" "3.844.83"
Set-Text "E3" "  +1.05%  "

# Row 4 - TetherUSD
Set-Text "E4" "  +0.00%  "

# Row 5 - BNB
Set-ForcedText "D5" "702.77"
Set-Text "E5" "  -0.82%  "

# Row 6 - Solana
Set-ForcedText "D6" "172.84"
Set-Text "E6" "  -0.33%  "

# Row 7 - LidoStakedEther
Set-Text "D7" "3.843.61"
Set-Text "E7" "  +1.05%  "

# Row 8 - USDC
Set-Text "E8" "  +0.02%  "

# Row 9 - XRP
Set-ForcedText "D9" "0.524"
Set-Text "E9" "  -1.18%  "

# Row 10 - Dogecoin
Set-Text "E10" "  -1.06%  "

# Row 11 - Toncoin
Set-ForcedText "D11" "7.35"
Set-Text "E11" "  -1.55%  "

# Row 12 - Cardano
Set-Text "E12" "  -0.89%  "

# Row 13 - ShibaInu
Set-ForcedText "D13" "0.0000256"
Set-Text "E13" "  -2.29%  "

# Row 14 - Avalanche
Set-ForcedText "D14" "36.44"
Set-Text "E14" "  +0.16%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-Text "D15" "4.496.13"
Set-Text "E15" "  +1.15%  "

# Row 16 - WrappedEther
Set-Text "D16" "3.762.41"
Set-Text "E16" "  -1.17%  "

# Row 17 - WrappedBTC
Set-Text "D17" "70.978.71"
Set-Text "E17" "  -0.12%  "

# Row 18 - Polkadot
Set-ForcedText "D18" "7.17"

# Row 19 - TRON
Set-Text "E19" "  +0.77%  "

# Row 20 - Chainlink
Set-ForcedText "D20" "17.35"
Set-Text "E20" "  -3.12%  "

# Row 21 - Uniswap
Set-ForcedText "D21" "10.68"
Set-Text "E21" "  -4.55%  "

# Row 22 - BitcoinCash
Set-ForcedText "D22" "492.25"
Set-Text "E22" "  +1.61%  "

# Row 23 - Polygon
Set-ForcedText "D23" "0.716"
Set-Text "E23" "  -0.17%  "

# Row 24 - Litecoin
Set-ForcedText "D24" "85.05"

# Row 25 - PEPE
Set-Text "E25" "  -0.35%  "

# Row 26 - InternetComputer(DFINITY)
Set-Text "E26" "  -2.47%  "

# Row 27 - RenderToken
Set-Text "E27" "  -0.18%  "

# Row 28 - Fetch.AI
Set-Text "E28" "  -3.14%  "

# Row 29 - PancakeSwap
Set-ForcedText "D29" "3.17"
Set-Text "E29" "  +1.79%  "

# Row 31 - NEARProtocol
Set-Text "E31" "  -1.44%  "

# Row 32 - ImmutableX
Set-ForcedText "D32" "2.28"
Set-Text "E32" "  -0.79%  "

# Row 33 - Kaspa
Set-Text "E33" "  +1.99%  "

# Row 34 - EthereumClassic
Set-ForcedText "D34" "29.41"
Set-Text "E34" "  -0.90%  "

# Row 35 - RenzoRestakedETH
Set-Text "D35" "3.803.66"
Set-Text "E35" "  +1.28%  "

# Row 36 - Aptos
Set-ForcedText "D36" "9.13"
Set-Text "E36" "  -1.48%  "

# Row 37 - Binance-PegBSC-USD
Set-Text "E37" "  -0.02%  "

# Row 38 - Hedera
Set-Text "E38" "  -0.41%  "

# Row 39 - Stacks
Set-Text "E39" "  +5.71%  "

# Row 40 - Filecoin
Set-Text "E40" "  +0.48%  "

# Row 41 - Mantle
Set-Text "E41" "  +5.76%  "

# Row 42 - dogwifhat
Set-Text "E42" "  -6.73%  "

# Row 43 - USDe
Set-Text "E43" "  +0.04%  "

# Row 44 - FirstDigitalUSD
Set-Text "E44" "  +0.07%  "

# Row 45 - Monero
Set-ForcedText "D45" "163.44"
Set-Text "E45" "  +0.50%  "

# Row 46 - FLOKI
Set-ForcedText "D46" "0.000309"
Set-Text "E46" "  -6.67%  "

# Row 47 - OKB
Set-ForcedText "D47" "48.72"
Set-Text "E47" "  -1.52%  "

# Row 48 - TheGraph
Set-Text "E48" "  -0.86%  "

# Row 49 - Cosmos
Set-ForcedText "D49" "8.62"
Set-Text "E49" "  +0.43%  "

# Row 50 - Bittensor
Set-ForcedText "D50" "412.23"
Set-Text "E50" "  +2.99%  "

# Row 51 - Arweave
Set-Text "E51" "  -4.14%  "
